# Auto-generated edit script: apply numeric updates from the commit diff
# to the Cactuar_Profits workbook (sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 300
$ws.Range("I18").Value = 300
$ws.Range("K18").Value = 300
$ws.Range("M18").Value = -16
$ws.Range("H19").Value = 1554.9
$ws.Range("I19").Value = 1988.4445
$ws.Range("J19").Value = 1200.1818
$ws.Range("K19").Value = 1988.4445
$ws.Range("L19").Value = 1200.1818
$ws.Range("M19").Value = -1813.4445
$ws.Range("N19").Value = -1550.1818
$ws.Range("H28").Value = 2790.8462
$ws.Range("I28").Value = 2742.4443
$ws.Range("J28").Value = 2899.75
$ws.Range("K28").Value = 2742.4443
$ws.Range("L28").Value = 2899.75
$ws.Range("M28").Value = -2257.4443
$ws.Range("N28").Value = -3869.75
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 2654.8333
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -2792.8333
$ws.Range("H136").Value = 99736.5
$ws.Range("J136").Value = 99736.5
$ws.Range("L136").Value = 99736.5
$ws.Range("N136").Value = -109936.5
$ws.Range("H137").Value = 199774.8
$ws.Range("I137").Value = 326052.53
$ws.Range("J137").Value = 4044.35
$ws.Range("K137").Value = 978157.5900000001
$ws.Range("L137").Value = 12133.05
$ws.Range("M137").Value = -975607.5900000001
$ws.Range("N137").Value = -17233.05
$ws.Range("H139").Value = 150000
$ws.Range("J139").Value = 150000
$ws.Range("L139").Value = 150000
$ws.Range("N139").Value = -160280
$ws.Range("H140").Value = 82529.89
$ws.Range("J140").Value = 84007.5
$ws.Range("L140").Value = 84007.5
$ws.Range("N140").Value = -94367.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 749
$ws.Range("I2").Value = 677.3333
$ws.Range("K2").Value = 677.3333
$ws.Range("M2").Value = -564.3333
$ws.Range("H32").Value = 16207.701
$ws.Range("I32").Value = 15211.804
$ws.Range("K32").Value = 15211.804
$ws.Range("M32").Value = -14924.804
$ws.Range("H61").Value = 5634
$ws.Range("I61").Value = 5327.2
$ws.Range("K61").Value = 5327.2
$ws.Range("M61").Value = -5115.2
$ws.Range("H97").Value = 16492.63
$ws.Range("I97").Value = 23461.154
$ws.Range("K97").Value = 23461.154
$ws.Range("M97").Value = -22965.154
$ws.Range("H101").Value = 72500
$ws.Range("J101").Value = 72500
$ws.Range("L101").Value = 72500
$ws.Range("N101").Value = -78990
$ws.Range("H102").Value = 1438.65
$ws.Range("I102").Value = 1399.2
$ws.Range("K102").Value = 1399.2
$ws.Range("M102").Value = 222.8
$ws.Range("H116").Value = 749
$ws.Range("I116").Value = 677.3333
$ws.Range("K116").Value = 677.3333
$ws.Range("M116").Value = 1616.6667
$ws.Range("H136").Value = 5634
$ws.Range("I136").Value = 5327.2
$ws.Range("K136").Value = 15981.6
$ws.Range("M136").Value = -13431.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 749
$ws.Range("I3").Value = 677.3333
$ws.Range("K3").Value = 677.3333
$ws.Range("M3").Value = -563.3333
$ws.Range("H26").Value = 7051.875
$ws.Range("I26").Value = 7051.875
$ws.Range("K26").Value = 7051.875
$ws.Range("M26").Value = -6759.875
$ws.Range("H105").Value = 2561.5557
$ws.Range("I105").Value = 2631.647
$ws.Range("K105").Value = 2631.647
$ws.Range("M105").Value = -884.6469999999999
$ws.Range("H134").Value = 5535.2246
$ws.Range("I134").Value = 2448.7727
$ws.Range("K134").Value = 7346.3181
$ws.Range("M134").Value = -4811.3181

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 5263
$ws.Range("I86").Value = 4877.125
$ws.Range("J86").Value = 5949
$ws.Range("K86").Value = 4877.125
$ws.Range("L86").Value = 5949
$ws.Range("M86").Value = -3754.125
$ws.Range("N86").Value = -8195
$ws.Range("H89").Value = 5263
$ws.Range("I89").Value = 4877.125
$ws.Range("J89").Value = 5949
$ws.Range("K89").Value = 24385.625
$ws.Range("L89").Value = 29745
$ws.Range("M89").Value = -18769.625
$ws.Range("N89").Value = -40977
$ws.Range("H99").Value = 3976.1282
$ws.Range("I99").Value = 2029.6154
$ws.Range("K99").Value = 2029.6154
$ws.Range("M99").Value = -531.6153999999999
$ws.Range("H122").Value = 4221.7827
$ws.Range("I122").Value = 2607.1333
$ws.Range("J122").Value = 7249.25
$ws.Range("K122").Value = 7821.3999
$ws.Range("L122").Value = 21747.75
$ws.Range("M122").Value = -5371.3999
$ws.Range("N122").Value = -26647.75
$ws.Range("H126").Value = 3976.1282
$ws.Range("I126").Value = 2029.6154
$ws.Range("K126").Value = 6088.8462
$ws.Range("M126").Value = -3618.8462
$ws.Range("H134").Value = 2943.818
$ws.Range("I134").Value = 2375.889
$ws.Range("K134").Value = 7127.667
$ws.Range("M134").Value = -4592.667
$ws.Range("H141").Value = 186345.05
$ws.Range("J141").Value = 196621
$ws.Range("L141").Value = 196621
$ws.Range("N141").Value = -206981

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 6067764.5
$ws.Range("I4").Value = 444268.75
$ws.Range("J4").Value = 8879512
$ws.Range("K4").Value = 1332806.25
$ws.Range("L4").Value = 26638536
$ws.Range("M4").Value = -1332694.25
$ws.Range("N4").Value = -26638760
$ws.Range("H75").Value = 62501700
$ws.Range("I75").Value = 799.5
$ws.Range("J75").Value = 83335330
$ws.Range("K75").Value = 2398.5
$ws.Range("L75").Value = 250005990
$ws.Range("M75").Value = -1400.5
$ws.Range("N75").Value = -250007986
$ws.Range("H78").Value = 62501700
$ws.Range("I78").Value = 799.5
$ws.Range("J78").Value = 83335330
$ws.Range("K78").Value = 7195.5
$ws.Range("L78").Value = 750017970
$ws.Range("M78").Value = -2203.5
$ws.Range("N78").Value = -750027954
$ws.Range("H87").Value = 22546.5
$ws.Range("I87").Value = 21093
$ws.Range("K87").Value = 63279
$ws.Range("M87").Value = -62031
$ws.Range("H90").Value = 22546.5
$ws.Range("I90").Value = 21093
$ws.Range("K90").Value = 189837
$ws.Range("M90").Value = -183597
$ws.Range("H93").Value = 6099.857
$ws.Range("J93").Value = 7000
$ws.Range("L93").Value = 21000
$ws.Range("N93").Value = -24744
$ws.Range("H117").Value = 4672.077
$ws.Range("I117").Value = 4076.5
$ws.Range("J117").Value = 4936.778
$ws.Range("K117").Value = 12229.5
$ws.Range("L117").Value = 14810.334
$ws.Range("M117").Value = -8787.5
$ws.Range("N117").Value = -21694.334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4849.857
$ws.Range("I70").Value = 4487.5
$ws.Range("J70").Value = 5333
$ws.Range("K70").Value = 4487.5
$ws.Range("L70").Value = 5333
$ws.Range("M70").Value = -4217.5
$ws.Range("N70").Value = -5873
$ws.Range("H73").Value = 4849.857
$ws.Range("I73").Value = 4487.5
$ws.Range("J73").Value = 5333
$ws.Range("K73").Value = 4487.5
$ws.Range("L73").Value = 5333
$ws.Range("M73").Value = -3551.5
$ws.Range("N73").Value = -7205
$ws.Range("H135").Value = 123852.664
$ws.Range("J135").Value = 123852.664
$ws.Range("L135").Value = 123852.664
$ws.Range("N135").Value = -133992.664
$ws.Range("H140").Value = 52031.5
$ws.Range("J140").Value = 63192.75
$ws.Range("L140").Value = 63192.75
$ws.Range("N140").Value = -73552.75
$ws.Range("H141").Value = 40619
$ws.Range("J141").Value = 40619
$ws.Range("L141").Value = 40619
$ws.Range("N141").Value = -50979

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 597.55554
$ws.Range("I55").Value = 195
$ws.Range("K55").Value = 195
$ws.Range("M55").Value = -22

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 166670580
$ws.Range("H138").Value = 98429
$ws.Range("J138").Value = 98429
$ws.Range("L138").Value = 98429
$ws.Range("N138").Value = -108709
$ws.Range("H140").Value = 98209.5
$ws.Range("J140").Value = 98209.5
$ws.Range("L140").Value = 98209.5
$ws.Range("N140").Value = -108569.5
